$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Status" column (C) cells that were marked "fail" are now marked "pass".
# Re-use the existing "pass" cell's number format / font / fill (copy format from
# C3, a passing cell) so the style entries that were unique to "fail" become
# unreferenced and get dropped from the workbook on save, then overwrite the
# cell text itself.
$passSource = $ws.Range("C3")
$failCells = @("C4","C6","C7","C9","C13","C16","C17","C18","C19","C20","C27","C28")

[void]$passSource.Copy()
foreach ($addr in $failCells) {
    $cell = $ws.Range($addr)
    [void]$cell.PasteSpecial(-4122)
    $cell.Value = "pass"
}

# Restore the active selection to where the user ended up after editing the sheet.
[void]$ws.Range("G27").Select()
